$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 20:22"

# Update the Ceuta row (row 55): Casos totales, Casos activos, Recuperados
$ws.Range("B55").Value = 95
$ws.Range("C55").Value = 21
$ws.Range("D55").Value = 70
